# Update "想去人数" (wish-to-go count) values on the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 70
        3  = 21695
        8  = 8086
        26 = 365
        28 = 75
        34 = 5179
        39 = 13360
        45 = 468
        48 = 335
    }
    "全部类型" = @{
        2  = 70
        3  = 21695
        7  = 8086
        23 = 365
        25 = 75
        34 = 5179
        39 = 13360
        45 = 468
        48 = 335
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $ws.Range("F$row").Value = $rowMap[$row]
    }
}
